$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right above the current row 33 (pushing the
# existing weekly records - and everything below them - down by two rows).
# This reproduces the diff exactly: rows 33..77 in the original file become
# rows 35..79, and the freshly inserted rows 33/34 hold this week's new
# "Murcott" price entries for Vega Monumental Concepción.
$ws.Range("A33:A34").EntireRow.Insert()

# New row 33: Murcott / Primera entry dated 2021-10-08 (serial 44477)
$ws.Range("A33").Value = 11
$ws.Range("B33").Value = "Vega Monumental Concepción"
$ws.Range("C33").Value = "Bíobío"
$ws.Range("D33").Value = 44477
$ws.Range("E33").Value = 8
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100102
$ws.Range("H33").Value = "Cítricos"
$ws.Range("I33").Value = 100102004
$ws.Range("J33").Value = "Mandarina"
$ws.Range("K33").Value = "Murcott"
$ws.Range("L33").Value = "Primera"
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 6000
$ws.Range("O33").Value = 6500
$ws.Range("P33").Value = 6250
$ws.Range("Q33").Value = "`$/bandeja 10 kilos"
$ws.Range("R33").Value = "Provincia de Limarí"
$ws.Range("S33").Value = 625
$ws.Range("T33").Value = 10

# New row 34: Murcott / Segunda entry dated 2021-10-08 (serial 44477)
$ws.Range("A34").Value = 11
$ws.Range("B34").Value = "Vega Monumental Concepción"
$ws.Range("C34").Value = "Bíobío"
$ws.Range("D34").Value = 44477
$ws.Range("E34").Value = 8
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100102
$ws.Range("H34").Value = "Cítricos"
$ws.Range("I34").Value = 100102004
$ws.Range("J34").Value = "Mandarina"
$ws.Range("K34").Value = "Murcott"
$ws.Range("L34").Value = "Segunda"
$ws.Range("M34").Value = 100
$ws.Range("N34").Value = 5000
$ws.Range("O34").Value = 5000
$ws.Range("P34").Value = 5000
$ws.Range("Q34").Value = "`$/bandeja 10 kilos"
$ws.Range("R34").Value = "Provincia de Limarí"
$ws.Range("S34").Value = 500
$ws.Range("T34").Value = 10
